$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the existing "ExpPoints" column (C), shifting
# it (and its header) from C to G, so we end up with:
#   A=Rank, B=Team, C=WIN, D=TOP2, E=TOP4, F=RELEGATION, G=ExpPoints
$ws.Range("C1:F1").EntireColumn.Insert()

# Set the new header labels (row 1). The inserted columns already inherit
# the bold/bordered header style from the column that used to be C.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP2"
$ws.Range("E1").Value = "TOP4"
$ws.Range("F1").Value = "RELEGATION"

# Make sure the new header cells carry the same style as the rest of the
# header row (bold font + border), matching columns A/B/G.
$ws.Range("C1:F1").Style = $ws.Range("A1").Style

# Reserve the new data cells (rows 2-19) for the upcoming Monte Carlo
# simulation results (WIN / TOP2 / TOP4 / RELEGATION probabilities). They
# stay blank for now, but keep them present as normal (unstyled) cells.
$ws.Range("C2:F19").Style = "Normal"
